$wb = $excel.ActiveWorkbook

# --- Sheet 1: Requirements Phase Defects ---
$ws1 = $wb.Worksheets.Item("Requirements Phase Defects")

$ws1.Range("C10").Value = "R01"
$ws1.Range("C11").Value = "R02"
$ws1.Range("C12").Value = "R03"
$ws1.Range("C13").Value = "R04"
$ws1.Range("C14").Value = "R05"
$ws1.Range("C15").Value = "R06"
$ws1.Range("C16").Value = "R07"
$ws1.Range("E10").Value = "Ce se intampla daca cele 8 mese sunt ocupate? - Mesaj de informare"
$ws1.Range("E13").Value = "De cate aplicatii este nevoie? Client, Chelner, Bucatar?"
$ws1.Range("E16").Value = "Meniul se incarca doar la pornirea aplicatiei? - da"

$ws1.Range("I3").Value = "Butacu Stefan"
$ws1.Range("J3").Value = 232
$ws1.Range("I4").Value = "Cadar Eduard"
$ws1.Range("J4").Value = 232

$ws1.Range("E26").Clear()

$ws1.Range("J4").Select()

# --- Sheet 2: Architect. Design Phase Defects ---
$ws2 = $wb.Worksheets.Item("Architect. Design Phase Defects")
$ws2.Range("I3").Select()

Write-Output "done"
